$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("measure_config")

# Add the rebase date and rebase comment for the second data row (row 2)
$ws.Range("L2").Value = '"2020-04-27"'
$ws.Range("M2").Value = "Rebased to demonstrate the method.  Add the rebase_dates and rebase_comment to 'measure_config.xlsx'."

# Column M needs to widen to fit the new, longer comment text
$ws.Columns.Item(13).ColumnWidth = 90
